$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.994.28"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.918.79"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.98"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4597"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3826"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9817"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.23"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.900.07"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.694"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.967"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06984"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.19"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009458"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.67"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").Value = "28.988.94"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.335"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "2.156.55"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.088"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.37"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.707"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.81"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.857"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09300"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8682"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.113"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.049"
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05710"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.154"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9998"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02041"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.061"
$ws.Range("E40").Value = "  +13.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.537"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5511"
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1752"
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000002998"
$ws.Range("E44").Value = "  +3.43%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.389"
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.204"
$ws.Range("E46").Value = "  +6.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5184"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.27"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06903"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.48"
$ws.Range("E51").Value = "  -0.14%  "
